$wb = $excel.ActiveWorkbook

# Sheet 1: "Cutting Speed" - add row 11 (B:E), all numeric
$ws1 = $wb.Worksheets.Item("Cutting Speed")
$ws1.Range("B11").Value = 210
$ws1.Range("C11").Value = 80
$ws1.Range("D11").Value = 6
$ws1.Range("E11").Value = 0.12

# Sheet 2: "Material Removal Rate" - add row 27 (B:F)
$ws2 = $wb.Worksheets.Item("Material Removal Rate")
$ws2.Range("B27").Value = 1
$ws2.Range("C27").Value = 60
$ws2.Range("D27").Value = 602
$ws2.Range("E27").Value = 36.12
$ws2.Range("F27").Value = "cm³/min"

# Sheet 3: "Helix Angle" - add row 12 (B:F)
$ws3 = $wb.Worksheets.Item("Helix Angle")
$ws3.Range("B12").Value = 5
$ws3.Range("C12").Value = 6
$ws3.Range("D12").Value = 0.06
$ws3.Range("E12").Value = 1.09
$ws3.Range("F12").Value = "°"

# Sheet 4: "Ramp Angle" - add rows 11 and 12 (B:E)
$ws4 = $wb.Worksheets.Item("Ramp Angle")
$ws4.Range("B11").Value = 1000
$ws4.Range("C11").Value = 5
$ws4.Range("D11").Value = 0.29
$ws4.Range("E11").Value = "°"

$ws4.Range("B12").Value = 250
$ws4.Range("C12").Value = 5
$ws4.Range("D12").Value = 1.15
$ws4.Range("E12").Value = "°"

# Sheet 5: "Surface Roughness" - add row 13 (B:D)
$ws5 = $wb.Worksheets.Item("Surface Roughness")
$ws5.Range("B13").Value = 0.2
$ws5.Range("C13").Value = 0.8
$ws5.Range("D13").Value = 2.08
